$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Date column (B) for rows 2 and 4 with new timestamp values,
# reflecting updated RAD Test Cases for Filing Year drop down and MD CRN changes.
$ws.Range("B2").Value = "Sat Feb 17 22:56:35 EST 2024"
$ws.Range("B4").Value = "Sat Feb 17 22:56:48 EST 2024"
